$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Acessar agência virtual'
$ws.Range('B2').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C2').Value = 'Seção ''O que é este serviço'' não encontrada'

$ws.Range('A3').Value = 'Acessar agência virtual'
$ws.Range('B3').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C3').Value = 'Seção ''Exigências'' não encontrada'

$ws.Range('A4').Value = 'Acessar agência virtual'
$ws.Range('B4').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C4').Value = 'Seção ''Quem pode utilizar'' não encontrada'

$ws.Range('A5').Value = 'Acessar agência virtual'
$ws.Range('B5').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C5').Value = 'Seção ''Prazos'' não encontrada'

$ws.Range('A6').Value = 'Acessar agência virtual'
$ws.Range('B6').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C6').Value = 'Seção ''Custos'' não encontrada'

$ws.Range('A7').Value = 'Acessar dados do portal da transparência'
$ws.Range('B7').Value = 'https://www.ms.gov.br/comunicacao-e-transparencia/acessar-dados-do-portal-da-transparencia171'
$ws.Range('C7').Value = 'Seção ''Outras Informações'' não encontrada'

$ws.Range('A8').Value = 'Acessar agência virtual'
$ws.Range('B8').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C8').Value = 'Seção ''Etapas'' não encontrada'

$ws.Range('A9').Value = 'Acessar agência virtual'
$ws.Range('B9').Value = 'https://www.ms.gov.br/energia/agencia-virtual22'
$ws.Range('C9').Value = 'Seção ''Outras Informações'' não encontrada'

$ws.Range('A10').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B10').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C10').Value = 'Seção ''O que é este serviço'' não encontrada'

$ws.Range('A11').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B11').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C11').Value = 'Seção ''Exigências'' não encontrada'

$ws.Range('A12').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B12').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C12').Value = 'Seção ''Quem pode utilizar'' não encontrada'

$ws.Range('A13').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B13').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C13').Value = 'Seção ''Prazos'' não encontrada'

$ws.Range('A14').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B14').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C14').Value = 'Seção ''Custos'' não encontrada'

$ws.Range('A15').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B15').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C15').Value = 'Seção ''Etapas'' não encontrada'

$ws.Range('A16').Value = 'Acessar gratuitamente programas da TV educativa (Portal da Educativa)'
$ws.Range('B16').Value = 'https://www.ms.gov.br/educacao-e-pesquisa/acessar-gratuitamente-programas-da-tv-educativa-portal-da-educativa175'
$ws.Range('C16').Value = 'Seção ''Outras Informações'' não encontrada'

Write-Host "done"